$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary section text + value updates ---
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 231739
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 7
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# --- Existing detail rows: update period/value data ---
$ws.Range("F16").Value = 35112
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2011"
# E19 stays "2010" (unchanged)

# --- Insert two new detail rows before the old last row (row 20) ---
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# Copy the formatting of a normal detail row (row 19) down into the two new rows
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new rows' values (same worker, new periods)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1050972722"
$ws.Range("D20").Value = "ANDRES FELIPE GAMARRA TINOCO"
$ws.Range("E20").Value = "2009"
$ws.Range("F20").Value = 35112
$ws.Range("G20").Value = 877803

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1050972722"
$ws.Range("D21").Value = "ANDRES FELIPE GAMARRA TINOCO"
$ws.Range("E21").Value = "2008"
$ws.Range("F21").Value = 35112
$ws.Range("G21").Value = 877803

# The old last row (now shifted to row 22) keeps its bordered style automatically;
# just update its period/value to the final one.
$ws.Range("E22").Value = "2007"
$ws.Range("F22").Value = 21067

# --- Footer (signature) rows, now at 27/28 after the insert ---
$ws.Range("B27").Value = "___________________________________"
$ws.Range("H27").Value = "___________________________________"
$ws.Range("B28").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H28").Value = "FIRMA DEL REPRESENTANTE LEGAL"
